# "sem ok para testar" - remove the "ok" values from column G (rows 2-4)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2:G4").ClearContents()
